$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.960.58"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.543.57"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'305.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").Value = "'98.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.39%  "
$ws.Range("D7").Value = "'0.578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'36.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").Value = "'0.0827"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.02%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'7.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.114"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "2.929.41"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "2.566.94"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "'15.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.19%  "
$ws.Range("D17").Value = "'0.874"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "42.947.19"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "'13.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").Value = "0.0₃0992"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").Value = "'6.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'71.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'254.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "'2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "'27.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").Value = "'38.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.23%  "
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").Value = "'157.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'2.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'19.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.35%  "
$ws.Range("D35").Value = "'0.0801"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").Value = "'3.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "'2.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.04%  "
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("D39").Value = "'24.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.13%  "
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  +20.51%  "
$ws.Range("D42").Value = "'3.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("D43").Value = "'3.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").Value = "2.093.24"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "'86.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("D48").Value = "'8.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").Value = "2.786.47"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'74.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.89%  "
$ws.Range("D51").Value = "'0.193"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.53%  "
